$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.047.24'
$ws.Range('E2').Value = '  -1.62%  '

$ws.Range('D3').Value = '1.554.67'
$ws.Range('E3').Value = '  -0.84%  '

$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('E5').Value = '  +0.12%  '

$ws.Range('D6').Value = '287.81'
$ws.Range('E6').Value = '  +0.29%  '

$ws.Range('D7').Value = '0.3963'
$ws.Range('E7').Value = '  +5.70%  '

$ws.Range('E8').Value = '  -1.78%  '

$ws.Range('D9').Value = '42.36'

$ws.Range('D10').Value = '1.111'
$ws.Range('E10').Value = '  -3.47%  '

$ws.Range('D11').Value = '0.07337'
$ws.Range('E11').Value = '  -1.22%  '

$ws.Range('D12').Value = '1.000'
$ws.Range('E12').Value = '  +0.03%  '

$ws.Range('D13').Value = '18.96'
$ws.Range('E13').Value = '  -7.53%  '

$ws.Range('D14').Value = '5.688'
$ws.Range('E14').Value = '  -2.79%  '

$ws.Range('D15').Value = '6.755'
$ws.Range('E15').Value = '  -1.23%  '

$ws.Range('D16').Value = '1.561.36'
$ws.Range('E16').Value = '  -0.62%  '

$ws.Range('D17').Value = '0.00001127'
$ws.Range('E17').Value = '  +2.36%  '

$ws.Range('D18').Value = '0.06609'
$ws.Range('E18').Value = '  -1.32%  '

$ws.Range('D19').Value = '85.04'
$ws.Range('E19').Value = '  -1.17%  '

$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = '1.000'
$ws.Range('E20').Value = '  +0.07%  '

$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '6.372'
$ws.Range('E21').Value = '  +0.16%  '

$ws.Range('E22').Value = '  -2.32%  '

$ws.Range('D23').Value = '11.32'
$ws.Range('E23').Value = '  -3.31%  '

$ws.Range('D24').Value = '22.047.59'

$ws.Range('D25').Value = '2.335'
$ws.Range('E25').Value = '  +0.71%  '

$ws.Range('D26').Value = '2.513'
$ws.Range('E26').Value = '  -2.22%  '

$ws.Range('D27').Value = '148.69'
$ws.Range('E27').Value = '  -2.20%  '

$ws.Range('D28').Value = '18.77'
$ws.Range('E28').Value = '  -2.95%  '

$ws.Range('D29').Value = '4.869'
$ws.Range('E29').Value = '  -1.03%  '

$ws.Range('D30').Value = '1.735.96'
$ws.Range('E30').Value = '  -0.64%  '

$ws.Range('D31').Value = '120.72'

$ws.Range('D32').Value = '1.079'
$ws.Range('E32').Value = '  +1.61%  '

$ws.Range('D33').Value = '5.752'
$ws.Range('E33').Value = '  -3.25%  '

$ws.Range('B34').Value = 'Stellar'
$ws.Range('C34').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D34').Value = '0.08337'
$ws.Range('E34').Value = '  +0.87%  '

$ws.Range('B35').Value = 'WEMIXTOKEN'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').Value = '1.654'
$ws.Range('E35').Value = '  -14.99%  '

$ws.Range('B36').Value = 'FraxShare'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D36').Value = '9.276'
$ws.Range('E36').Value = '  -4.09%  '

$ws.Range('D37').Value = '0.06205'
$ws.Range('E37').Value = '  -2.01%  '

$ws.Range('D38').Value = '0.02279'
$ws.Range('E38').Value = '  -4.62%  '

$ws.Range('D39').Value = '5.153'
$ws.Range('E39').Value = '  -2.48%  '

$ws.Range('D40').Value = '0.2101'
$ws.Range('E40').Value = '  -4.28%  '

$ws.Range('D41').Value = '1.213'
$ws.Range('E41').Value = '  -6.01%  '

$ws.Range('D42').Value = '1.001'
$ws.Range('E42').Value = '  -0.04%  '

$ws.Range('D43').Value = '10.71'
$ws.Range('E43').Value = '  -3.72%  '

$ws.Range('D44').Value = '0.5885'
$ws.Range('E44').Value = '  -3.75%  '

$ws.Range('D45').Value = '13.33'
$ws.Range('E45').Value = '  -2.69%  '

$ws.Range('D46').Value = '3.723'
$ws.Range('E46').Value = '  -0.72%  '

$ws.Range('D47').Value = '0.5646'
$ws.Range('E47').Value = '  -4.72%  '

$ws.Range('D48').Value = '1.919'
$ws.Range('E48').Value = '  -4.66%  '

$ws.Range('D49').Value = '117.95'
$ws.Range('E49').Value = '  -4.93%  '

$ws.Range('D50').Value = '1.147'
$ws.Range('E50').Value = '  -2.90%  '

$ws.Range('D51').Value = '0.06854'
$ws.Range('E51').Value = '  -4.22%  '
